$wb = $excel.ActiveWorkbook

# --- Rename the data sheet 30.12.20 -> 31.12.20 (this also updates the
# defined name "Bundesländer001" which refers to '30.12.20'!$A$1:$G$17) ---
$wsExpl = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)
$wsData.Name = "31.12.20"

# --- "Erläuterung" sheet updates ---
# C6: time of data stand moved from 08:30 to 12:30
$wsExpl.Range("C6").Value = "12:30 Uhr"
# A10: updated wording of the "Achtung" footnote (adds "oder Korrekturen")
$wsExpl.Range("A10").Value = "Achtung: Die Differenz zum Vortag kann Nachmeldungen oder Korrekturen aus vorangegangenen Tagen enthalten und spiegelt nicht immer die innerhalb des Vortags tatsächlich durchgeführte Zahl der Impfungen wider. "

# --- Data sheet (31.12.20) updates: new day's figures ---
function Set-Row($ws, $row, $b, $c, $d, $e, $f, $g) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    if ($null -ne $d) { $ws.Cells.Item($row, 4).Value = $d }
    $ws.Cells.Item($row, 5).Value = $e
    if ($null -ne $f) { $ws.Cells.Item($row, 6).Value = $f }
    $ws.Cells.Item($row, 7).Value = $g
}

Set-Row $wsData 2  17086 4295 7578  5423  1037 2693
Set-Row $wsData 3  37955 9749 7840  18450 718  14265
Set-Row $wsData 4  11114 2030 7029  2952  105  8161
Set-Row $wsData 5  3219  183  204   2995  20   224
Set-Row $wsData 6  1741  50   821   408   6    237
Set-Row $wsData 7  2759  719  1129  1420  $null 1341
Set-Row $wsData 8  21373 5699 6073  12345 821  8513
Set-Row $wsData 9  11494 393  178   6583  8    4804
Set-Row $wsData 10 3566  861  567   1981  654  2044
Set-Row $wsData 11 24924 3033 $null 8963  $null 16334
Set-Row $wsData 12 5112  865  $null 2329  $null 2783
Set-Row $wsData 13 2716  597  2065  248   $null 1080
Set-Row $wsData 14 3290  637  226   2754  1    535
Set-Row $wsData 15 11146 1046 3506  5582  458  5346
Set-Row $wsData 16 7270  1689 1766  4523  1527 2817

$wsData.Cells.Item(17, 3).Value = 0
$wsData.Cells.Item(17, 6).Value = 0

# New remark cell about missing KV-Nordrhein data for NRW row
$wsData.Cells.Item(11, 8).Value = "(für den 31.12. wurden keine Daten aus der KV-Nordrhein übermittelt)"
